$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, $Value)
    $Range.NumberFormat = "@"
    $Range.Value = $Value
}

# Row 13 and 14 swap Coin name and Link (Chainlink <-> WrappedEther)
Set-TextValue $ws.Range("B13") "Chainlink"
Set-TextValue $ws.Range("C13") "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue $ws.Range("B14") "WrappedEther"
Set-TextValue $ws.Range("C14") "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"

Set-TextValue $ws.Range("D2") "29.650.86"
Set-TextValue $ws.Range("E2") "  -2.35%  "
Set-TextValue $ws.Range("D3") "2.009.81"
Set-TextValue $ws.Range("E3") "  -4.66%  "
Set-TextValue $ws.Range("D4") "1.012"
Set-TextValue $ws.Range("E4") "  +0.44%  "
Set-TextValue $ws.Range("D5") "332.05"
Set-TextValue $ws.Range("E6") "  +0.46%  "
Set-TextValue $ws.Range("D7") "0.5047"
Set-TextValue $ws.Range("E7") "  -3.42%  "
Set-TextValue $ws.Range("D8") "0.4273"
Set-TextValue $ws.Range("E8") "  -3.71%  "
Set-TextValue $ws.Range("D9") "54.85"
Set-TextValue $ws.Range("E9") "  +0.49%  "
Set-TextValue $ws.Range("D10") "0.09215"
Set-TextValue $ws.Range("E10") "  -1.91%  "
Set-TextValue $ws.Range("D11") "1.129"
Set-TextValue $ws.Range("E11") "  -3.77%  "
Set-TextValue $ws.Range("D12") "23.68"
Set-TextValue $ws.Range("E12") "  -5.14%  "
Set-TextValue $ws.Range("D13") "8.159"
Set-TextValue $ws.Range("E13") "  -6.04%  "
Set-TextValue $ws.Range("D14") "2.001.00"
Set-TextValue $ws.Range("E14") "  -1.78%  "
Set-TextValue $ws.Range("D15") "6.563"
Set-TextValue $ws.Range("E15") "  -5.49%  "
Set-TextValue $ws.Range("D16") "95.55"
Set-TextValue $ws.Range("E16") "  -6.18%  "
Set-TextValue $ws.Range("E17") "  +0.42%  "
Set-TextValue $ws.Range("D18") "0.00001129"
Set-TextValue $ws.Range("E18") "  -2.94%  "
Set-TextValue $ws.Range("D19") "0.06676"
Set-TextValue $ws.Range("E19") "  -0.80%  "
Set-TextValue $ws.Range("D20") "19.96"
Set-TextValue $ws.Range("E20") "  -5.94%  "
Set-TextValue $ws.Range("D21") "1.010"
Set-TextValue $ws.Range("E21") "  +0.48%  "
Set-TextValue $ws.Range("D22") "5.999"
Set-TextValue $ws.Range("E22") "  -5.56%  "
Set-TextValue $ws.Range("D23") "29.643.42"
Set-TextValue $ws.Range("E23") "  -2.47%  "
Set-TextValue $ws.Range("D24") "12.10"
Set-TextValue $ws.Range("E24") "  -4.28%  "
Set-TextValue $ws.Range("E25") "  -1.27%  "
Set-TextValue $ws.Range("D26") "159.43"
Set-TextValue $ws.Range("E26") "  -2.32%  "
Set-TextValue $ws.Range("D27") "20.86"
Set-TextValue $ws.Range("E27") "  -5.30%  "
Set-TextValue $ws.Range("D28") "6.446"
Set-TextValue $ws.Range("E28") "  -5.35%  "
Set-TextValue $ws.Range("D29") "2.339"
Set-TextValue $ws.Range("E29") "  -7.74%  "
Set-TextValue $ws.Range("D30") "129.05"
Set-TextValue $ws.Range("E30") "  -3.68%  "
Set-TextValue $ws.Range("D31") "1.069"
Set-TextValue $ws.Range("E31") "  -7.25%  "
Set-TextValue $ws.Range("D32") "1.586"
Set-TextValue $ws.Range("E32") "  -8.64%  "
Set-TextValue $ws.Range("D33") "0.09981"
Set-TextValue $ws.Range("E33") "  -5.33%  "
Set-TextValue $ws.Range("D34") "5.878"
Set-TextValue $ws.Range("E34") "  -6.22%  "
Set-TextValue $ws.Range("D35") "3.816"
Set-TextValue $ws.Range("E35") "  -2.83%  "
Set-TextValue $ws.Range("D36") "9.572"
Set-TextValue $ws.Range("E36") "  -8.12%  "
Set-TextValue $ws.Range("D37") "0.02481"
Set-TextValue $ws.Range("E37") "  -5.52%  "
Set-TextValue $ws.Range("D38") "1.321"
Set-TextValue $ws.Range("E38") "  -2.46%  "
Set-TextValue $ws.Range("D39") "0.06404"
Set-TextValue $ws.Range("E39") "  -5.65%  "
Set-TextValue $ws.Range("D40") "0.6613"
Set-TextValue $ws.Range("E40") "  -6.16%  "
Set-TextValue $ws.Range("D41") "11.82"
Set-TextValue $ws.Range("E41") "  -5.99%  "
Set-TextValue $ws.Range("D42") "0.2079"
Set-TextValue $ws.Range("E42") "  -6.58%  "
Set-TextValue $ws.Range("D43") "1.010"
Set-TextValue $ws.Range("E43") "  +0.47%  "
Set-TextValue $ws.Range("D44") "0.6394"
Set-TextValue $ws.Range("E44") "  -6.66%  "
Set-TextValue $ws.Range("D45") "13.55"
Set-TextValue $ws.Range("E45") "  -5.98%  "
Set-TextValue $ws.Range("D46") "2.224"
Set-TextValue $ws.Range("E46") "  -5.76%  "
Set-TextValue $ws.Range("D47") "1.299"
Set-TextValue $ws.Range("E47") "  -4.56%  "
Set-TextValue $ws.Range("D48") "3.532"
Set-TextValue $ws.Range("E48") "  -3.14%  "
Set-TextValue $ws.Range("D49") "0.07011"
Set-TextValue $ws.Range("E49") "  -3.33%  "
Set-TextValue $ws.Range("D50") "0.00000000323"
Set-TextValue $ws.Range("E50") "  -7.73%  "
Set-TextValue $ws.Range("D51") "1.136"
Set-TextValue $ws.Range("E51") "  -5.83%  "
